$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header / summary block updates
# ---------------------------------------------------------------------------
# Valor Mora total (row 11)
$ws.Range("E11").Value = 496342

# Cant. Trabajadores (row 13, C13) and Cant. Periodos (row 13, F13)
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 12

# ---------------------------------------------------------------------------
# 2) Insert a new data row at the top of the table (row 16) for the new
#    worker FARID SIR ARRIETA, pushing the existing rows down by one.
#    Copy formatting from the row that will land at 17 (the old row 16,
#    which carries the regular interior-row style) so the new row matches.
# ---------------------------------------------------------------------------
$ws.Rows.Item(16).Insert()
$ws.Range("B17:J17").Copy()
$ws.Range("B16:J16").PasteSpecial(-4122)

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1129565103"
$ws.Range("D16").Value = "FARID SIR ARRIETA"
$ws.Range("E16").Value = "1811"
$ws.Range("F16").Value = 192000
$ws.Range("G16").Value = 1100000

# ---------------------------------------------------------------------------
# 3) Re-populate rows 17-27 with YOIBER MIRANDA GUERRA's periods, now listed
#    in descending order (1712 down to 1702). Styling for these rows is
#    already correct because it simply shifted down with the insert above.
# ---------------------------------------------------------------------------
$yoiberPeriods = @("1712","1711","1710","1709","1708","1707","1706","1705","1704","1703","1702")
for ($i = 0; $i -lt $yoiberPeriods.Length; $i++) {
    $r = 17 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "1143353859"
    $ws.Range("D$r").Value = "YOIBER MIRANDA GUERRA"
    $ws.Range("E$r").Value = $yoiberPeriods[$i]
    $ws.Range("F$r").Value = 27578
    $ws.Range("G$r").Value = 737717
}

# ---------------------------------------------------------------------------
# 4) Row 28 already carries the special "last row" bordered style (it was
#    the old row 27 before the insert shifted it down). Populate it with
#    EDGARDO JOSE FERNANDEZ SALGUEDO's entry, which moves to the bottom.
# ---------------------------------------------------------------------------
$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "73083976"
$ws.Range("D28").Value = "EDGARDO JOSE FERNANDEZ SALGUEDO"
$ws.Range("E28").Value = "1706"
$ws.Range("F28").Value = 984
$ws.Range("G28").Value = 737717

Write-Host "Edit applied"
